$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column by column to preserve shared string ordering
# Column A
$ws.Range("A1").Value = "username"
$ws.Range("A2").Value = "Manish"

# Column B
$ws.Range("B1").Value = "password"
$ws.Range("B2").Value = "kk"
$ws.Range("B3").Value = "PP"

# Column C
$ws.Range("C1").Value = "is_correct"
$ws.Range("C2").Value = "Y"
$ws.Range("C3").Value = "N"

# Column D
$ws.Range("D1").Value = "age"
$ws.Range("D2").Value = 22
$ws.Range("D3").Value = 23

# Column E
$ws.Range("E1").Value = "gender"
$ws.Range("E2").Value = "M"
$ws.Range("E3").Value = "F"

# A3 last, since it's a new/unique string added after gender column values
$ws.Range("A3").Value = "TTT"

# Column widths to match bestFit (AutoFit columns A and B)
$ws.Columns.Item(1).ColumnWidth = 9.0
$ws.Columns.Item(2).ColumnWidth = 17.9

# Set selection to A3
$ws.Range("A3").Select()
